# #3473 replaced two properties that had gaps
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell-level data corrections (rows 2-10) ---
$ws.Range("C2").Value = "Medstar POB South Tower"
$ws.Range("I2").Value = 1985
$ws.Range("L2").Value = 76319

$ws.Range("C3").Value = "1801 Pennsylvania Avenue, LLC"

$ws.Range("E4").Value = "300 E ST SW"
$ws.Range("J4").Value = "TWO INDEPENDENCE HANA OW LLC"
$ws.Range("L4").Value = 627655

$ws.Range("L5").Value = 58717

$ws.Range("C6").Value = "Hampton House"
$ws.Range("E6").Value = "2700 CONNECTICUT AVENUE NW"
$ws.Range("H6").Value = 20008
$ws.Range("I6").Value = 1921
$ws.Range("J6").Value = "2700 CONECTICUT AVENUE LLC"
$ws.Range("L6").Value = 83580

$ws.Range("H7").Value = 20005
$ws.Range("L7").Value = 145697

$ws.Range("E8").Value = "1428 H ST NW"

# Row 10 replaced wholesale: "DPW Vehicle Maintenance Facility 2" -> "School Without Walls @ Francis Stevens"
$ws.Range("C10").Value = "School Without Walls @ Francis Stevens"
$ws.Range("D10").Value = "K-12 School"
$ws.Range("E10").Value = "2425 N STREET NW"
$ws.Range("H10").Value = 20037
$ws.Range("I10").Value = 1924
$ws.Range("J10").Value = "DISTRICT OF COLUMBIA"
$ws.Range("L10").Value = 127991

# --- Remove the stray date-ish formatting on the Year Built column ---
$ws.Range("I2:I10").ClearFormats()

# --- Drop the now-unused "EUI Target Year" column (M) entirely ---
$ws.Range("M1").EntireColumn.Delete()

# --- Selection / view state ---
$ws.Range("A1:L10").Select()
